# Exp4/Dados/Buck/2.xlsx - "Inicio relatorio - Experimento 4"
# Convert the E-column text labels into real numbers, then add a new
# "Ia" (current) column with a reference value and a ratio formula
# F = E / $I$1 formatted to three decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E3:E6 used to hold the numbers as text (shared strings); replace
#     them with real numeric values so they can be used in formulas.
$ws.Range("E3").Value = 0.62
$ws.Range("E4").Value = 0.76
$ws.Range("E5").Value = 0.9
$ws.Range("E6").Value = 1.02

# --- New header / reference cells.
$ws.Range("F1").Value = "Ia"
$ws.Range("H1").Value = "R:"
$ws.Range("I1").Value = 5.3

# --- Ratio formulas in column F, formatted with three decimal places.
$ws.Range("F3").Formula = "=E3/`$I`$1"
$ws.Range("F4:F6").Formula = "=E4/`$I`$1"

$ws.Range("F3:F6").NumberFormat = "0.000"

# --- Match the selection left behind by the author.
$ws.Range("F3:F6").Select()
